$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Golang Developer", "https://www.dice.com/job-detail/1722b428-adee-4657-abdd-99d4ec26f95a", "Remote", "Contract, Third Party", "70+", "Ace Technologies, Inc."),
    @("Go Developer", "https://www.dice.com/job-detail/0872e306-f80e-43ea-9821-966922024476", "Alpharetta, Georgia", "Contract", "60+", "AgreeYa Solutions"),
    @("Software Engineer (Backend-First, React Optional)", "https://www.dice.com/job-detail/d8489b2a-b2f0-422b-8a0c-a3cd99c669e7", "Fremont, California", "Contract, Third Party", "$60 - $70", "e-IT Professionals Corp.")
)

$startRow = 118
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
